$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supplier-portal document submission run: the automation stamps the
# latest "User_On_<date>-<time>" login/submission marker into A2 (the
# sheet's single status cell), joining the existing history of
# User_On_* shared strings already recorded in this test-data sheet.
$ws.Range("A2").Value = "User_On_03/01/19-11:24"
